$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, center/top alignment) from N1 to O1:P1
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1:P1").PasteSpecial(-4122) | Out-Null

# Header values
$ws.Cells.Item(1, 15).Value = 14
$ws.Cells.Item(1, 16).Value = 15

# Data values for columns O (15) and P (16), rows 2-67
$ws.Cells.Item(2, 15).Value = -0.2125358924799186
$ws.Cells.Item(2, 16).Value = -0.2115618058472211
$ws.Cells.Item(3, 15).Value = 0.2804199035622875
$ws.Cells.Item(3, 16).Value = 0.280266812961422
$ws.Cells.Item(4, 15).Value = 0.2471288954660506
$ws.Cells.Item(4, 16).Value = 0.2493461651665066
$ws.Cells.Item(5, 15).Value = -0.03066674603132093
$ws.Cells.Item(5, 16).Value = -0.03215856623565759
$ws.Cells.Item(6, 15).Value = 0.2292678284262772
$ws.Cells.Item(6, 16).Value = 0.2285342342747512
$ws.Cells.Item(7, 15).Value = -0.3708471285936795
$ws.Cells.Item(7, 16).Value = -0.3695668025220744
$ws.Cells.Item(8, 15).Value = -0.1414754737630154
$ws.Cells.Item(8, 16).Value = -0.1370980462934937
$ws.Cells.Item(9, 15).Value = -0.2824211976526718
$ws.Cells.Item(9, 16).Value = -0.2792724227708773
$ws.Cells.Item(10, 15).Value = 0.4296922439413407
$ws.Cells.Item(10, 16).Value = 0.4284785154670471
$ws.Cells.Item(11, 15).Value = -0.1771803635111997
$ws.Cells.Item(11, 16).Value = -0.1768736292806702
$ws.Cells.Item(12, 15).Value = -0.007993258584414957
$ws.Cells.Item(12, 16).Value = -0.01311420331645243
$ws.Cells.Item(13, 15).Value = -0.02012032951476822
$ws.Cells.Item(13, 16).Value = -0.02240596973031777
$ws.Cells.Item(14, 15).Value = 0.2034617189598873
$ws.Cells.Item(14, 16).Value = 0.1949363713022942
$ws.Cells.Item(15, 15).Value = 0.1075175228514393
$ws.Cells.Item(15, 16).Value = 0.09365039824090399
$ws.Cells.Item(16, 15).Value = 0.5568604849942538
$ws.Cells.Item(16, 16).Value = 0.5423456070667526
$ws.Cells.Item(17, 15).Value = 0.6353618880560403
$ws.Cells.Item(17, 16).Value = 0.6187708266261263
$ws.Cells.Item(18, 15).Value = -0.06822152723612694
$ws.Cells.Item(18, 16).Value = -0.07649961690124031
$ws.Cells.Item(19, 15).Value = 0.4111767471347207
$ws.Cells.Item(19, 16).Value = 0.4024112839341616
$ws.Cells.Item(20, 15).Value = 0.4582147822908585
$ws.Cells.Item(20, 16).Value = 0.4414293750204749
$ws.Cells.Item(21, 15).Value = 0.6719310804902803
$ws.Cells.Item(21, 16).Value = 0.6585691956035908
$ws.Cells.Item(22, 15).Value = 0.4420989245126605
$ws.Cells.Item(22, 16).Value = 0.4278287041609329
$ws.Cells.Item(23, 15).Value = -0.01163862315642603
$ws.Cells.Item(23, 16).Value = -0.02414307508162968
$ws.Cells.Item(24, 15).Value = 2.127493328782096
$ws.Cells.Item(24, 16).Value = 1.932502094650812
$ws.Cells.Item(25, 15).Value = 0.2968006167949152
$ws.Cells.Item(25, 16).Value = 0.2940663808138195
$ws.Cells.Item(26, 15).Value = 0.1562846535522882
$ws.Cells.Item(26, 16).Value = 0.1475281067981714
$ws.Cells.Item(27, 15).Value = 0.05468087244996441
$ws.Cells.Item(27, 16).Value = 0.04558693761122762
$ws.Cells.Item(28, 15).Value = 0.7957114587545457
$ws.Cells.Item(28, 16).Value = 0.7881270009708706
$ws.Cells.Item(29, 15).Value = 1.944803223047363
$ws.Cells.Item(29, 16).Value = 1.813507612086966
$ws.Cells.Item(30, 15).Value = 0.6340952467610315
$ws.Cells.Item(30, 16).Value = 0.6271246598854809
$ws.Cells.Item(31, 15).Value = -0.4988264367412496
$ws.Cells.Item(31, 16).Value = -0.4994610635860754
$ws.Cells.Item(32, 15).Value = 0.5298711885311274
$ws.Cells.Item(32, 16).Value = 0.5242255491853003
$ws.Cells.Item(33, 15).Value = 0.7386313586358176
$ws.Cells.Item(33, 16).Value = 0.7359939609253274
$ws.Cells.Item(34, 15).Value = -0.8530485429394171
$ws.Cells.Item(34, 16).Value = -0.8569084167357776
$ws.Cells.Item(35, 15).Value = 0.7747917777023493
$ws.Cells.Item(35, 16).Value = 0.7767229460071331
$ws.Cells.Item(36, 15).Value = 0.7136182366490784
$ws.Cells.Item(36, 16).Value = 0.7180691500632012
$ws.Cells.Item(37, 15).Value = 0.6766867498219962
$ws.Cells.Item(37, 16).Value = 0.6810933295711001
$ws.Cells.Item(38, 15).Value = 0.624868485521312
$ws.Cells.Item(38, 16).Value = 0.6251823317021289
$ws.Cells.Item(39, 15).Value = 0.5771589433185644
$ws.Cells.Item(39, 16).Value = 0.5795754927826959
$ws.Cells.Item(40, 15).Value = 0.7395180674214268
$ws.Cells.Item(40, 16).Value = 0.7415396913555018
$ws.Cells.Item(41, 15).Value = 0.5501742733505088
$ws.Cells.Item(41, 16).Value = 0.5531472627949815
$ws.Cells.Item(42, 15).Value = 0.5711413899992752
$ws.Cells.Item(42, 16).Value = 0.573955310691755
$ws.Cells.Item(43, 15).Value = 0.6508280297409779
$ws.Cells.Item(43, 16).Value = 0.6524120255491189
$ws.Cells.Item(44, 15).Value = 0.6651797703601553
$ws.Cells.Item(44, 16).Value = 0.6687613046237939
$ws.Cells.Item(45, 15).Value = 0.6229599254771434
$ws.Cells.Item(45, 16).Value = 0.6303016033025318
$ws.Cells.Item(46, 15).Value = -1.289877892817069
$ws.Cells.Item(46, 16).Value = -1.292050892517576
$ws.Cells.Item(47, 15).Value = -1.007611343809836
$ws.Cells.Item(47, 16).Value = -1.009208081669638
$ws.Cells.Item(48, 15).Value = -0.8817417231319653
$ws.Cells.Item(48, 16).Value = -0.8812152176844328
$ws.Cells.Item(49, 15).Value = -0.6434093130117455
$ws.Cells.Item(49, 16).Value = -0.6432022674863436
$ws.Cells.Item(50, 15).Value = -0.06197223061920823
$ws.Cells.Item(50, 16).Value = -0.06334908601762014
$ws.Cells.Item(51, 15).Value = -0.8689872467926884
$ws.Cells.Item(51, 16).Value = -0.8678481333333029
$ws.Cells.Item(52, 15).Value = -0.8689872467926884
$ws.Cells.Item(52, 16).Value = -0.8678481333333029
$ws.Cells.Item(53, 15).Value = -1.143470861639018
$ws.Cells.Item(53, 16).Value = -1.143916963618337
$ws.Cells.Item(54, 15).Value = -0.1762935326014256
$ws.Cells.Item(54, 16).Value = -0.1763967384019015
$ws.Cells.Item(55, 15).Value = -1.026735395816189
$ws.Cells.Item(55, 16).Value = -1.028668341483162
$ws.Cells.Item(56, 15).Value = -0.8955526673001041
$ws.Cells.Item(56, 16).Value = -0.8985876543085083
$ws.Cells.Item(57, 15).Value = -0.9037702065448942
$ws.Cells.Item(57, 16).Value = -0.9093909733198395
$ws.Cells.Item(58, 15).Value = -1.048734878463608
$ws.Cells.Item(58, 16).Value = -1.053493814648325
$ws.Cells.Item(59, 15).Value = -0.7808261565040457
$ws.Cells.Item(59, 16).Value = -0.7813981875539462
$ws.Cells.Item(60, 15).Value = -0.4036403583600968
$ws.Cells.Item(60, 16).Value = -0.4053190666441998
$ws.Cells.Item(61, 15).Value = 0.3876347372667817
$ws.Cells.Item(61, 16).Value = 0.3875143449590815
$ws.Cells.Item(62, 15).Value = -1.154432337354456
$ws.Cells.Item(62, 16).Value = -1.160843475630652
$ws.Cells.Item(63, 15).Value = -0.5620656377423728
$ws.Cells.Item(63, 16).Value = -0.5582259242165966
$ws.Cells.Item(64, 15).Value = -0.8375786247702771
$ws.Cells.Item(64, 16).Value = -0.8377355649355009
$ws.Cells.Item(65, 15).Value = -0.02216534937065239
$ws.Cells.Item(65, 16).Value = -0.02326897064305818
$ws.Cells.Item(66, 15).Value = -0.7320954311482674
$ws.Cells.Item(66, 16).Value = -0.7386429067919765
$ws.Cells.Item(67, 15).Value = -0.7013653139109522
$ws.Cells.Item(67, 16).Value = -0.710394028805811
